$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 68; this shifts existing rows 68-101 down to 69-102
# (the sheet's used range grows from A1:R101 to A1:R102, matching the target diff).
$ws.Rows.Item(68).Insert()

# Populate the newly inserted row 68 with the new weekly record.
$ws.Cells.Item(68, 1).Value = 5
$ws.Cells.Item(68, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(68, 3).Value = "Maule"
$ws.Cells.Item(68, 4).Value = 44606
$ws.Cells.Item(68, 5).Value = 7
$ws.Cells.Item(68, 6).Value = 100112001
$ws.Cells.Item(68, 7).Value = "Berenjena"
$ws.Cells.Item(68, 8).Value = "Sin especificar"
$ws.Cells.Item(68, 9).Value = "Primera"
$ws.Cells.Item(68, 10).Value = 200
$ws.Cells.Item(68, 11).Value = 7000
$ws.Cells.Item(68, 12).Value = 7000
$ws.Cells.Item(68, 13).Value = 7000
$ws.Cells.Item(68, 14).Value = '$/caja 60 unidades'
$ws.Cells.Item(68, 15).Value = 'Región del Maule'
$ws.Cells.Item(68, 16).Value = 117
$ws.Cells.Item(68, 17).Value = 60
$ws.Cells.Item(68, 18).Value = "Hortaliza"
